$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "flower/flower002.jpg"
$ws.Range("C2").Value = "haken"
$ws.Range("D2").Value = "flower"
$ws.Range("B3").Value = "flower/flower012.jpg"
$ws.Range("C3").Value = "liefern"
$ws.Range("D3").Value = "flower"
$ws.Range("B4").Value = "flower/flower020.jpg"
$ws.Range("C4").Value = "scheitern"
$ws.Range("D4").Value = "flower"
$ws.Range("B5").Value = "face/face024.jpg"
$ws.Range("C5").Value = "tagen"
$ws.Range("D5").Value = "face"
$ws.Range("B6").Value = "face/face015.jpg"
$ws.Range("C6").Value = "stärken"
$ws.Range("D6").Value = "face"
$ws.Range("B7").Value = "face/face005.jpg"
$ws.Range("C7").Value = "schicken"
$ws.Range("D7").Value = "face"
$ws.Range("B8").Value = "face/face022.jpg"
$ws.Range("C8").Value = "pflegen"
$ws.Range("D8").Value = "face"
$ws.Range("B9").Value = "flower/flower014.jpg"
$ws.Range("C9").Value = "gründen"
$ws.Range("D9").Value = "flower"
$ws.Range("B10").Value = "face/face020.jpg"
$ws.Range("C10").Value = "segeln"
$ws.Range("D10").Value = "face"
$ws.Range("B11").Value = "face/face013.jpg"
$ws.Range("C11").Value = "gelten"
$ws.Range("D11").Value = "face"
$ws.Range("B12").Value = "flower/flower019.jpg"
$ws.Range("C12").Value = "regnen"
$ws.Range("D12").Value = "flower"
$ws.Range("B13").Value = "face/face012.jpg"
$ws.Range("C13").Value = "jubeln"
$ws.Range("D13").Value = "face"
$ws.Range("B14").Value = "flower/flower004.jpg"
$ws.Range("C14").Value = "kaufen"
$ws.Range("D14").Value = "flower"
$ws.Range("B15").Value = "face/face026.jpg"
$ws.Range("C15").Value = "krachen"
$ws.Range("D15").Value = "face"
$ws.Range("B16").Value = "face/face002.jpg"
$ws.Range("C16").Value = "bitten"
$ws.Range("D16").Value = "face"
$ws.Range("B17").Value = "face/face000.jpg"
$ws.Range("C17").Value = "starten"
$ws.Range("D17").Value = "face"
$ws.Range("B18").Value = "face/face025.jpg"
$ws.Range("C18").Value = "opfern"
$ws.Range("D18").Value = "face"
$ws.Range("B19").Value = "flower/flower028.jpg"
$ws.Range("C19").Value = "töten"
$ws.Range("D19").Value = "flower"
$ws.Range("B20").Value = "flower/flower001.jpg"
$ws.Range("C20").Value = "fühlen"
$ws.Range("D20").Value = "flower"
$ws.Range("B21").Value = "face/face001.jpg"
$ws.Range("C21").Value = "ehren"
$ws.Range("D21").Value = "face"
$ws.Range("B22").Value = "flower/flower007.jpg"
$ws.Range("C22").Value = "stechen"
$ws.Range("D22").Value = "flower"
$ws.Range("B23").Value = "flower/flower010.jpg"
$ws.Range("C23").Value = "saufen"
$ws.Range("D23").Value = "flower"
$ws.Range("B24").Value = "flower/flower018.jpg"
$ws.Range("C24").Value = "sieben"
$ws.Range("D24").Value = "flower"
$ws.Range("B25").Value = "face/face023.jpg"
$ws.Range("C25").Value = "kehren"
$ws.Range("D25").Value = "face"
$ws.Range("B26").Value = "flower/flower015.jpg"
$ws.Range("C26").Value = "spielen"
$ws.Range("D26").Value = "flower"
$ws.Range("B27").Value = "flower/flower029.jpg"
$ws.Range("C27").Value = "rücken"
$ws.Range("D27").Value = "flower"
$ws.Range("B28").Value = "face/face009.jpg"
$ws.Range("C28").Value = "wiegen"
$ws.Range("D28").Value = "face"
$ws.Range("B29").Value = "flower/flower006.jpg"
$ws.Range("C29").Value = "langen"
$ws.Range("D29").Value = "flower"
$ws.Range("B30").Value = "face/face017.jpg"
$ws.Range("C30").Value = "hupen"
$ws.Range("D30").Value = "face"
$ws.Range("B31").Value = "flower/flower021.jpg"
$ws.Range("C31").Value = "hoffen"
$ws.Range("D31").Value = "flower"
$ws.Range("B32").Value = "face/face003.jpg"
$ws.Range("C32").Value = "husten"
$ws.Range("D32").Value = "face"
$ws.Range("B33").Value = "flower/flower022.jpg"
$ws.Range("C33").Value = "währen"
$ws.Range("D33").Value = "flower"
